# Refactor profit computation and electricity price extraction
# Adds a new "sTimePeriods" worksheet (hourly electricity price lookup)
# right before the "Coordinates" sheet.

$wb = $excel.ActiveWorkbook

$coords = $wb.Worksheets.Item("Coordinates")
$chargingStations = $wb.Worksheets.Item("sChargingStations")
$pathTypes = $wb.Worksheets.Item("sPathTypes")

# Insert the new sheet immediately before "Coordinates" so it lands in the
# same slot Coordinates used to occupy (rId6), pushing Coordinates to rId7.
$new = $wb.Worksheets.Add($coords)
$new.Name = "sTimePeriods"

# --- Headers -------------------------------------------------------------
$new.Range("A1").Value = "pPeriod"
$new.Range("B1").Value = "aux_pElectricityCost [€/MWh]"
$new.Range("C1").Value = "pElectricityCost [$/kWh]"

# Match existing header styling used elsewhere in the workbook:
#  - bold + themed text color (same as the index/id columns, e.g. sChargingStations!A1)
#  - bold, un-themed text (same as sChargingStations!E1:F1)
$chargingStations.Range("A1").Copy() | Out-Null
$new.Range("A1").PasteSpecial(-4122) | Out-Null
$new.Range("A1").Value = "pPeriod"

$chargingStations.Range("E1").Copy() | Out-Null
$new.Range("B1:C1").PasteSpecial(-4122) | Out-Null
$new.Range("B1").Value = "aux_pElectricityCost [€/MWh]"
$new.Range("C1").Value = "pElectricityCost [$/kWh]"

# --- Data rows (hour-of-day 0-23 electricity price lookup) ---------------
$periods = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23)
$prices  = @(19.12,28.68,23.6,22.85,24.24,26.15,27.95,36.91,49.89,53.1,28.39,13.28,12.23,11.97,12.46,10.58,9.27,9.93,8.21,10.36,13.76,43.83,48.02,32)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = $i + 2
    $new.Cells.Item($r, 1).Value = $periods[$i]
    $new.Cells.Item($r, 2).Value = $prices[$i]
    $new.Cells.Item($r, 3).Formula = "=B$r*0.00142"
}

# Style column A (period number) like other id/index columns.
$chargingStations.Range("A1").Copy() | Out-Null
$new.Range("A2:A25").PasteSpecial(-4122) | Out-Null
for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = $i + 2
    $new.Cells.Item($r, 1).Value = $periods[$i]
}

# Style column B (raw €/MWh price) like other plain numeric columns (2 decimals).
$pathTypes.Range("H2").Copy() | Out-Null
$new.Range("B2:B25").PasteSpecial(-4122) | Out-Null
for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = $i + 2
    $new.Cells.Item($r, 2).Value = $prices[$i]
}

# Style column C (converted $/kWh) with the same base font, then apply the
# new 3-decimal number format used specifically for this column.
$pathTypes.Range("H2").Copy() | Out-Null
$new.Range("C2:C25").PasteSpecial(-4122) | Out-Null
$new.Range("C2:C25").NumberFormat = "0.000"
for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = $i + 2
    $new.Cells.Item($r, 3).Formula = "=B$r*0.00142"
}

# --- Column widths (best-fit, matching authored widths) ------------------
$new.Columns.Item(1).ColumnWidth = 18.33203125
$new.Columns.Item(2).ColumnWidth = 26
$new.Columns.Item(3).ColumnWidth = 21.77734375

# --- Sheet view: zoom + selection, and make this the active/visible tab --
$new.Activate()
$excel.ActiveWindow.Zoom = 115
$new.Range("E19").Select()

$wb.Application.CutCopyMode = $false
